$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Mã phiếu" (code) from P2210-00003 to P2210-00001
$ws.Range("A2").Value = "P2210-00001"

# Update "Nội dung phiếu" (content) - now mirrors the code value instead of
# the old "Phiếu nhập tháng 10" text
$ws.Range("C2").Value = "P2210-00001"

# "Ngày nhập" (import date) stays the same text value
$ws.Range("D2").Value = "11-10-2022 00:00:00"

# Update "Tổng tiền nhập" (total amount) from 1.100.000 VND to 1.200.000 VND
$ws.Range("E2").Value = "1.200.000 VND"

# Clear "Ghi chú" (note) - previously duplicated the content text
$ws.Range("F2").Value = ""

# Adjust column widths to match new content (the engine adds ~5/7 padding
# chars on top of the value we set, so back that out to land on the exact
# target widths of 17 and 9)
$ws.Columns.Item(3).ColumnWidth = 16.285714285714285
$ws.Columns.Item(6).ColumnWidth = 8.285714285714286
